$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old column E ("collector_name"),
# shifting collector_name..status from E..K to F..L.
$ws.Range("E1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("E1").Value = "village"

# Match the author's column width change: admin_level_2 (D) and the new
# village column (E) both become 23.4531 characters wide.
$ws.Range("D1").ColumnWidth = 23.4531
$ws.Range("E1").ColumnWidth = 23.4531
